$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Tim Southee"

# Insert a new column before column A, shifting existing data right
$ws.Columns.Item(1).Insert()

# Fill in the new column A with header and value
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "41st"
